# Applies the "Updated symbol list" edit: 98 cell updates across the
# coin ranking table (prices, 1h volume %, and several swapped
# coin name/link rows) taken from the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h %) store numbers-as-text in this sheet
# (e.g. "332.19", "1.19%"). Excel normally auto-converts such literals to a
# real number/percentage on assignment, so force the cell to Text format
# first, assign the literal, then restore the default "Normal" style so the
# cell formatting matches the original (unstyled) cells.
function Set-TextValue($ref, $value) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue 'D2' '332.19'
Set-TextValue 'E2' '1.19%'
Set-TextValue 'D3' '45.90'
Set-TextValue 'E3' '4.40%'
Set-TextValue 'D4' '5.639'
Set-TextValue 'E4' '2.31%'
Set-TextValue 'D5' '0.08364'
Set-TextValue 'E5' '4.31%'
Set-TextValue 'D6' '2.057'
Set-TextValue 'E6' '3.30%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.483'
Set-TextValue 'E7' '3.40%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D8' '0.9835'
Set-TextValue 'E8' '3.60%'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D9' '2.561'
Set-TextValue 'E9' '-2.00%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1157'
Set-TextValue 'E10' '2.26%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1933'
Set-TextValue 'E11' '3.41%'
$ws.Range('B12').Value = 'MCDex'
$ws.Range('C12').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D12' '10.40'
Set-TextValue 'E12' '-2.57%'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D13' '0.09962'
Set-TextValue 'E13' '1.58%'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D14' '0.04664'
Set-TextValue 'E14' '-0.81%'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D15' '0.1058'
Set-TextValue 'E15' '-0.58%'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D16' '0.001293'
Set-TextValue 'E16' '2.24%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D17' '0.006107'
Set-TextValue 'E17' '3.21%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D18' '3.374'
Set-TextValue 'E18' '0.50%'
Set-TextValue 'E19' '-3.16%'
Set-TextValue 'E20' '-0.32%'
Set-TextValue 'D21' '0.2653'
Set-TextValue 'E21' '4.26%'
Set-TextValue 'D22' '0.04201'
Set-TextValue 'D23' '0.001313'
Set-TextValue 'E23' '4.31%'
Set-TextValue 'D24' '0.004637'
Set-TextValue 'E24' '7.22%'
Set-TextValue 'D26' '0.0003749'
Set-TextValue 'E26' '0.17%'
Set-TextValue 'D38' '0.02777'
Set-TextValue 'E38' '7.47%'
Set-TextValue 'E39' '2.82%'
Set-TextValue 'D40' '0.007752'
Set-TextValue 'E40' '2.79%'
Set-TextValue 'D41' '0.1436'
Set-TextValue 'E41' '2.78%'
Set-TextValue 'D42' '0.007257'
Set-TextValue 'E42' '-3.80%'
Set-TextValue 'D43' '0.002015'
Set-TextValue 'E43' '0.03%'
Set-TextValue 'D44' '0.008104'
Set-TextValue 'E44' '-5.61%'
Set-TextValue 'D46' '0.00007314'
Set-TextValue 'E46' '2.24%'
Set-TextValue 'D47' '0.00000000752'
Set-TextValue 'E47' '0.28%'
Set-TextValue 'D48' '0.0005814'
Set-TextValue 'E48' '0.04%'
$ws.Range('B49').Value = 'CoinbaseStockToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue 'D49' '0.003506'
Set-TextValue 'E49' '-0.65%'
$ws.Range('B50').Value = 'BOLO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue 'D50' '0.003501'
Set-TextValue 'E50' '-2.99%'
Set-TextValue 'D51' '0.00002105'
Set-TextValue 'E51' '0.28%'

Write-Output "Applied 98 cell updates"
